$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Harish@566#"
$ws.Range("B3").Value = "Swetha@633#"
$ws.Range("B4").Value = "Swetha@202#"
$ws.Range("B5").Value = "Deepa@184#"
$ws.Range("B6").Value = "Rahul@466#"
$ws.Range("B7").Value = "Lakshmi@108#"
$ws.Range("B8").Value = "Sandeep@977#"
$ws.Range("B9").Value = "Preeti@540#"
$ws.Range("B10").Value = "Ajay@522#"
$ws.Range("B11").Value = "Rina@054#"
$ws.Range("B12").Value = "Vivek@854#"
$ws.Range("B13").Value = "Meena@989#"
$ws.Range("B14").Value = "Sahil@783#"
$ws.Range("B15").Value = "Tanvi@342#"
$ws.Range("B16").Value = "Yash@967#"
$ws.Range("B17").Value = "Anita@536#"
$ws.Range("B18").Value = "Rohit@836#"
$ws.Range("B19").Value = "Kiran@397#"
$ws.Range("B20").Value = "Sunita@926#"
$ws.Range("B21").Value = "Amitabh@917#"
$ws.Range("B22").Value = "Priya@227#"
$ws.Range("B23").Value = "Nitin@276#"
$ws.Range("B24").Value = "Sneha@137#"
$ws.Range("B25").Value = "Arjun@085#"
$ws.Range("B26").Value = "Divya@819#"
